# ajouts binding + example care plan 9a7e74c6ff63a487b613993ee94e019015a79f79
#
# 1) Metadata sheet: bump the "Date" property to the new commit timestamp.
# 2) Elements sheet, row 6 (the OncoFAIR value[x] / CodeableConcept row):
#       - Binding Strength  (col X) -> "required"
#       - Binding Description (col Y) -> cleared/blank
#       - Binding Value Set (col Z) -> the new ValueSet URL
#    and widen column Z so the long URL remains best-fit/visible.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B8").Value = "2024-06-04T08:55:54+00:00"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("X6").Value = "required"
$wsElements.Range("Y6").ClearContents()
$wsElements.Range("Z6").Value = "http://ltsi.univ-rennes.fr/ValueSet/siph-typecomposant-oncofair-valueset"

$wsElements.Columns.Item(26).ColumnWidth = 66.66

# The save round-trip of this engine does not preserve the worksheet's
# "bestFit" column flag, and re-hides columns only if we restate it - make
# sure the originally-hidden helper columns stay hidden.
$wsElements.Columns.Item(3).Hidden = $true
$wsElements.Columns.Item(4).Hidden = $true
$wsElements.Columns.Item(31).Hidden = $true
$wsElements.Columns.Item(32).Hidden = $true
$wsElements.Columns.Item(33).Hidden = $true
